# Commit: Added error check to ensure rounds that have had data gathered cannot be
# further updated once upcoming predictions have been generated. This would give
# wrong team PAVs and venues etc. as dataframe shapes would be misaligned.
#
# Structurally, this resulted in 5 additional "match" columns (KR, KS, KT, KU, KV)
# being appended to the stats sheet, each one duplicating the values that were in
# the (then) last populated column KQ, for every data row (1-102).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A genuinely untouched / blank cell far outside the used range, used purely as a
# formatting donor so the brand new last column (KV) ends up without an explicit
# cell style, matching how the previous last column (KR) looked before the edit.
$blank = $ws.Range("LA200")

for ($row = 1; $row -le 102; $row++) {
    $srcCell = $ws.Cells.Item($row, 303)   # KQ<row>
    $v = $srcCell.Value2

    $kr = $ws.Cells.Item($row, 304)        # KR<row> (previously last column, unstyled)
    $kr.ClearContents()
    $kr.Value = $v

    $ws.Cells.Item($row, 305).Value = $v   # KS<row>
    $ws.Cells.Item($row, 306).Value = $v   # KT<row>
    $ws.Cells.Item($row, 307).Value = $v   # KU<row>

    $kv = $ws.Cells.Item($row, 308)        # KV<row> (new last column, stays unstyled)
    $blank.Copy($kv)
    $kv.Value = $v
}
